# RPAR_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclosure note (A18)
#    from 2021-05-26 to 2021-05-27
#  - refresh the Weight (D) / Percent Change (E) model holdings figures
#    in rows 2-15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; lift protection so the cells can be written,
# then restore protection once the edits are in place.
$ws.Unprotect()

$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.0567623126348397
$ws.Range("E2").Value = 0.00184060371801964
$ws.Range("D3").Value = 0.02384620484897514
$ws.Range("E3").Value = 0.001916810427448956
$ws.Range("D4").Value = 0.03144102421967067
$ws.Range("E4").Value = 0.001869158878504695
$ws.Range("D5").Value = 0.03231928335934787
$ws.Range("E5").Value = 0.0007677543186179747
$ws.Range("D6").Value = 0.03613856790059093
$ws.Range("E6").Value = 0.02580225498699051
$ws.Range("D7").Value = 0.01872683558870062
$ws.Range("E7").Value = 0.002440810349036093
$ws.Range("D8").Value = 0.004472265314237267
$ws.Range("E8").Value = -0.0008760402978537085
$ws.Range("D9").Value = 0.006892865039726762
$ws.Range("E9").Value = 0.006062902614626831
$ws.Range("D10").Value = 0.07385212542564387
$ws.Range("E10").Value = -0.0005305039787797394
$ws.Range("D11").Value = 0.07396966196478283
$ws.Range("E11").Value = -0.0005296610169490679
$ws.Range("D12").Value = 0.1454527731460037
$ws.Range("E12").Value = -0.004058324953311265
$ws.Range("D13").Value = 0.3816744446036121
$ws.Range("E13").Value = -0.00174504842509382
$ws.Range("D14").Value = 0.1144516359538684
$ws.Range("E14").Value = -0.0004145470055568756
$ws.Range("D15").Value = 0.9999999999999998
$ws.Range("E15").Value = -0.0001323330834550607

$ws.Protect()
